$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FirstName / MiddleName(Initial) / LastName / EmailAddress for each
# registrant row. Password column (E) and header row (1) are unchanged.
$ws.Range("A2").Value = "Sam"
$ws.Range("B2").Value = "Ron"
$ws.Range("C2").Value = "Jin"
$ws.Range("D2").Value = "a120@email.com"

$ws.Range("A3").Value = "Peter"
$ws.Range("B3").Value = "Kio"
$ws.Range("C3").Value = "Conery"
$ws.Range("D3").Value = "a121@email.com"

$ws.Range("A4").Value = "Butna"
$ws.Range("B4").Value = "Amy"
$ws.Range("C4").Value = "Swan"
$ws.Range("D4").Value = "a122@email.com"

$ws.Range("A5").Value = "Jen"
$ws.Range("B5").Value = "Loper"
$ws.Range("C5").Value = "Kou"
$ws.Range("D5").Value = "a123@email.com"

# Match the author's final on-screen selection/view state as closely as the
# COM surface allows (scroll position past the data into row 8).
$ws.Range("E8").Select()
